# Add the new inventory item "STLS4V" (Extractor de abolladuras) as row 45,
# extending the sheet's used range from A1:J44 to A1:J45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45").Value = "STLS4V"
$ws.Range("B45").Value = "Extractor de abolladuras"
# C45 (Descripción) is intentionally left blank, matching row 44's pattern.
$ws.Range("D45").Value = 22000
$ws.Range("E45").Value = 75000
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = 22
$ws.Range("H45").Formula = "=(E45-D45)*G45"
$ws.Range("I45").Formula = "=D45*F45"
$ws.Range("J45").Value = 22000
